$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 294, pushing existing rows 294+ down to 296+
$ws.Rows.Item(294).Resize(2).Insert()

# Row 294 (Primera) - new data point
$ws.Cells.Item(294, 1).Value = 8
$ws.Cells.Item(294, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(294, 3).Value = "Coquimbo"
$ws.Cells.Item(294, 4).Value = 44889
$ws.Cells.Item(294, 5).Value = 4
$ws.Cells.Item(294, 6).Value = 100114014
$ws.Cells.Item(294, 7).Value = "Betarraga"
$ws.Cells.Item(294, 8).Value = "Sin especificar"
$ws.Cells.Item(294, 9).Value = "Primera"
$ws.Cells.Item(294, 10).Value = 2000
$ws.Cells.Item(294, 11).Value = 550
$ws.Cells.Item(294, 12).Value = 600
$ws.Cells.Item(294, 13).Value = 575
$ws.Cells.Item(294, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(294, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(294, 16).Value = 192
$ws.Cells.Item(294, 17).Value = 3
$ws.Cells.Item(294, 18).Value = "Hortaliza"

# Row 295 (Segunda) - new data point
$ws.Cells.Item(295, 1).Value = 8
$ws.Cells.Item(295, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(295, 3).Value = "Coquimbo"
$ws.Cells.Item(295, 4).Value = 44889
$ws.Cells.Item(295, 5).Value = 4
$ws.Cells.Item(295, 6).Value = 100114014
$ws.Cells.Item(295, 7).Value = "Betarraga"
$ws.Cells.Item(295, 8).Value = "Sin especificar"
$ws.Cells.Item(295, 9).Value = "Segunda"
$ws.Cells.Item(295, 10).Value = 1500
$ws.Cells.Item(295, 11).Value = 450
$ws.Cells.Item(295, 12).Value = 500
$ws.Cells.Item(295, 13).Value = 475
$ws.Cells.Item(295, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(295, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(295, 16).Value = 158
$ws.Cells.Item(295, 17).Value = 3
$ws.Cells.Item(295, 18).Value = "Hortaliza"

# Apply the same date number format (style) used in column D to the two new D cells
$ws.Range("D296").Copy()
$ws.Range("D294:D295").PasteSpecial(-4122) # xlPasteFormats
